# Applies per-cell value updates to cryptos.xlsx worksheet
# to reflect refreshed crypto price/volume data, plus a couple of
# row re-ordering / coin swaps (Stellar <-> EthereumClassic; Algorand -> Cronos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.331.10"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.933.80"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7551"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3172"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06998"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7766"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08019"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "1.932.62"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.339"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").Value = "30.347.92"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007921"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.744"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "2.181.48"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.660"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.471"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1333"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.195"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.371"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.515"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.390"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.116"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05140"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.282"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7500"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.772"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.804"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "77.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.424"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4456"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.963"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8323"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.785"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.478"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "975.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.38%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06041"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.80%  "

Write-Host "Applied 101 cell updates to cryptos sheet"
